$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2520
$ws1.Range("F7").Value = 53
$ws1.Range("F8").Value = 48
$ws1.Range("F12").Value = 614
$ws1.Range("F13").Value = 1479
$ws1.Range("F14").Value = 1261
$ws1.Range("F15").Value = 3
$ws1.Range("G15").Value = 1
$ws1.Range("F16").Value = 497
$ws1.Range("F17").Value = 3595
$ws1.Range("F18").Value = 642
$ws1.Range("F19").Value = 3284
$ws1.Range("F21").Value = 2140
$ws1.Range("F23").Value = 288
$ws1.Range("F25").Value = 3
$ws1.Range("G25").Value = 1
$ws1.Range("F26").Value = 1129
$ws1.Range("F29").Value = 979
$ws1.Range("F30").Value = 959
$ws1.Range("F31").Value = 73

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 211
$ws2.Range("F19").Value = 238

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 502

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 502
$ws4.Range("F12").Value = 2520
$ws4.Range("F13").Value = 2520
$ws4.Range("F15").Value = 53
$ws4.Range("F16").Value = 48
$ws4.Range("F27").Value = 1479
$ws4.Range("F29").Value = 1261
$ws4.Range("F30").Value = 497
$ws4.Range("F32").Value = 3595
$ws4.Range("F33").Value = 642
$ws4.Range("F34").Value = 3284
$ws4.Range("F36").Value = 2140
$ws4.Range("F38").Value = 288
$ws4.Range("F40").Value = 1129
$ws4.Range("F43").Value = 238
$ws4.Range("F48").Value = 979
$ws4.Range("F49").Value = 959
$ws4.Range("F50").Value = 73
